# Apply the "Structure of the Model Hierarchy" diagram rework:
#  - retitle to (V2) / add a second subtitle line
#  - reflow the three top-level rounded rectangles (taller / narrower, shifted up)
#  - reposition the "Site/Treatment/Tree Level" labels and the
#    "PAR, VPD" / "Rain Event" / "Sap Flow" callouts to match the new layout
#  - replace the single "Soil Moisture" left-arrow with three separate callouts:
#    "Soil Moisture (shallow)", "Soil Moisture (deep)" and "Pre-rain Totals"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$EMU = 12700.0

# --- Title textbox: reposition/resize + new two-paragraph text -------------
$title = $s.Shapes.Item("TextBox 6")
$title.Left = 1193800 / $EMU
$title.Top = 270933 / $EMU
$title.Width = 6874933 / $EMU

$tr = $title.TextFrame.TextRange
$tr.Text = "Structure of the Model Hierarchy ("
$tr.InsertAfter("V2)") | Out-Null
$title.TextFrame.TextRange.InsertAfter([char]13 + "Once for JUMO, Once for PIED") | Out-Null
$title.TextFrame.TextRange.InsertAfter(" ") | Out-Null

# --- Three top-level rounded rectangles: taller, narrower, shifted up ------
$r3 = $s.Shapes.Item("Rounded Rectangle 3")
$r3.Left = 524934 / $EMU
$r3.Top = 1413933 / $EMU
$r3.Width = 2201334 / $EMU
$r3.Height = 3327400 / $EMU

$r4 = $s.Shapes.Item("Rounded Rectangle 4")
$r4.Left = 3458634 / $EMU
$r4.Top = 1413933 / $EMU
$r4.Width = 2311399 / $EMU
$r4.Height = 3327400 / $EMU

$r5 = $s.Shapes.Item("Rounded Rectangle 5")
$r5.Left = 6498166 / $EMU
$r5.Top = 1413933 / $EMU
$r5.Width = 2294466 / $EMU
$r5.Height = 3327400 / $EMU

# --- "Site Level" / "Treatment Level" / "Tree Level" labels ---------------
$siteLevel = $s.Shapes.Item("TextBox 7")
$siteLevel.Left = 6993466 / $EMU
$siteLevel.Top = 1667933 / $EMU

$treatmentLevel = $s.Shapes.Item("TextBox 8")
$treatmentLevel.Left = 3598334 / $EMU
$treatmentLevel.Top = 1679601 / $EMU

$treeLevel = $s.Shapes.Item("TextBox 9")
$treeLevel.Left = 135468 / $EMU
$treeLevel.Top = 1679601 / $EMU

# --- "PAR, VPD" / "Rain Event" / "Sap Flow" callouts ------------------------
$parVpd = $s.Shapes.Item("Left Arrow 11")
$parVpd.Left = 6688665 / $EMU
$parVpd.Top = 2175932 / $EMU

$rainEvent = $s.Shapes.Item("Left Arrow 12")
$rainEvent.Left = 3598334 / $EMU
$rainEvent.Top = 2060601 / $EMU

$sapFlow = $s.Shapes.Item("Rounded Rectangle 13")
$sapFlow.Left = 1193800 / $EMU
$sapFlow.Top = 2904066 / $EMU

# --- Replace the single "Soil Moisture" arrow with three callouts ---------
$soilMoisture = $s.Shapes.Item("Left Arrow 18")

$shallow = $soilMoisture.Duplicate().Item(1)
$shallow.Left = 3598334 / $EMU
$shallow.Top = 2904066 / $EMU
$shallowTr = $shallow.TextFrame.TextRange
$shallowTr.Text = "Soil "
$shallowTr.InsertAfter("Moisture (shallow)") | Out-Null

$deep = $soilMoisture.Duplicate().Item(1)
$deep.Left = 3598334 / $EMU
$deep.Top = 3759200 / $EMU
$deepTr = $deep.TextFrame.TextRange
$deepTr.Text = "Soil "
$deepTr.InsertAfter("Moisture (deep)") | Out-Null

$preRain = $soilMoisture.Duplicate().Item(1)
$preRain.Left = 6688665 / $EMU
$preRain.Top = 3268133 / $EMU
$preRain.TextFrame.TextRange.Text = "Pre-rain Totals"

$soilMoisture.Delete()
